$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing row 12 timestamp (A12) with the refreshed value from
# the scheduled task run.
$ws.Range("A12").Value = 45864.66703143519

# Append the new scheduled-task reading as row 13.
$ws.Range("A13").Value = 45864.70865358998
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat

$ws.Range("B13").Value = 2025
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 19.94
$ws.Range("E13").Value = 70.8
$ws.Range("F13").Value = 116.02
$ws.Range("G13").Value = 12.87
$ws.Range("H13").Value = "ESE"
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "17:00:27"
